{"js": "// The document contains a single-column table. Three header-ish cells get\n// their text changed, ten brand-new rows are inserted right after the third\n// row, and the three trailing \"tab separated\" summary cells collapse down to\n// a single plain value each.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// 1) First three cells: \"100\" -> \"0M\", \"0\" -> \"0M\", \"256\" -> \"0M\".\ntable.getCell(0, 0).value = \"0M\";\ntable.getCell(1, 0).value = \"0M\";\ntable.getCell(2, 0).value = \"0M\";\nawait context.sync();\n\n// 2) Insert ten new rows right after the (now updated) third row, each\n//    holding one of the new measurement values.\nconst newValues = [\n  [\"133\"],\n  [\"0.00002\"],\n  [\"0.00005\"],\n  [\"0.00004\"],\n  [\"0.00001\"],\n  [\"0.00003\"],\n  [\"0.00004\"],\n  [\"0.00004\"],\n  [\"0.00488\"],\n  [\"100.0\"],\n];\nrows.items[2].insertRows(\"After\", newValues.length, newValues);\nawait context.sync();\n\n// 3) The three trailing rows (originally tab-separated multi-run cells)\n//    collapse to a single simple value each. After the ten-row insert they\n//    are now the last three rows of the table.\nconst tables2 = context.document.body.tables;\ntables2.load(\"items\");\nawait context.sync();\nconst table2 = tables2.items[0];\nconst rows2 = table2.rows;\nrows2.load(\"items\");\nawait context.sync();\n\nconst n = rows2.items.length;\ntable2.getCell(n - 3, 0).value = \"100\";\ntable2.getCell(n - 2, 0).value = \"0\";\ntable2.getCell(n - 1, 0).value = \"256\";\nawait context.sync();\n", "ps1": "# The document holds a single-column table. The first three cells get new\n# text, ten brand-new rows are inserted right after the (updated) third row,\n# and the three trailing \"tab separated\" summary cells collapse down to a\n# single plain value each.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# 1) First three cells: \"100\" -> \"0M\", \"0\" -> \"0M\", \"256\" -> \"0M\".\n$t.Cell(1,1).Range.Text = \"0M\"\n$t.Cell(2,1).Range.Text = \"0M\"\n$t.Cell(3,1).Range.Text = \"0M\"\n\n# 2) Insert ten new rows right after the (now updated) third row, each\n#    holding one of the new measurement values.\n$newValues = @(\"133\",\"0.00002\",\"0.00005\",\"0.00004\",\"0.00001\",\"0.00003\",\"0.00004\",\"0.00004\",\"0.00488\",\"100.0\")\n\n$insertBeforeIndex = 4\nforeach ($v in $newValues) {\n    $refRow = $t.Rows.Item($insertBeforeIndex)\n    $newRow = $t.Rows.Add($refRow)\n    $newRow.Cells.Item(1).Range.Text = $v\n    $insertBeforeIndex = $insertBeforeIndex + 1\n}\n\n# 3) The three trailing rows (originally tab-separated multi-run cells)\n#    collapse to a single simple value each. After the ten-row insert they\n#    are now the last three rows of the table.\n$rowCount = $t.Rows.Count\n$t.Cell($rowCount - 2, 1).Range.Text = \"100\"\n$t.Cell($rowCount - 1, 1).Range.Text = \"0\"\n$t.Cell($rowCount, 1).Range.Text = \"256\"\n"}
